$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for new "Save" column (copy formatting from neighboring header cell)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Save flag values per row (2..16)
$saveValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 1
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
